$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for the new "Save" column in H1, matching style of existing headers (s="1")
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate H2:H46 with Save flag: 1 if sum (column G) is >= 8, else 0
for ($r = 2; $r -le 46; $r++) {
    $g = $ws.Cells.Item($r, 7).Value2
    if ($g -ge 8) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
